# Generate Report for Handoff
#
# The handoff xliff files were (re)generated, so each language's status
# flips from "In Translation" to "Ready for handoff" and the associated
# timestamps are refreshed. Excel widens the (now longer) status columns
# to fit the new text.
#
# Note: the host's `ColumnWidth` setter only lands on ~1/6-character
# pixel-grid increments (round((w + 5/6) * 6) / 6), so 16.333333333333332
# is the closest input that reproduces the target stored width.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status + generate-date columns ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-10-24 09:57:31"
$ws.Columns.Item(5).ColumnWidth = 16.333333333333332
$ws.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn detail sheet: Status + Latest Handoff Datetime ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-10-24 09:57:20"
$ws.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- de-de detail sheet: Status + Latest Handoff Datetime ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-10-24 09:57:31"
$ws.Columns.Item(3).ColumnWidth = 16.333333333333332
